$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1477.125
$ws.Range("J19").Value = 1324
$ws.Range("L19").Value = 1324
$ws.Range("N19").Value = -1674

$ws.Range("H51").Value = 5833.1665
$ws.Range("I51").Value = 4999.75
$ws.Range("J51").Value = 7500
$ws.Range("K51").Value = 4999.75
$ws.Range("L51").Value = 7500
$ws.Range("M51").Value = -4515.75
$ws.Range("N51").Value = -8468

$ws.Range("H54").Value = 6500
$ws.Range("I54").Value = 6500
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 6500
$ws.Range("L54").Value = 0
$ws.Range("N54").Value = -6014
$ws.Range("M54").ClearContents()

$ws.Range("H86").Value = 2659.6
$ws.Range("J86").Value = 1600
$ws.Range("L86").Value = 1600
$ws.Range("N86").Value = -3846

$ws.Range("H89").Value = 2659.6
$ws.Range("J89").Value = 1600
$ws.Range("L89").Value = 8000
$ws.Range("N89").Value = -19232

$ws.Range("H106").Value = 2199
$ws.Range("I106").Value = 3994.5
$ws.Range("J106").Value = 1800
$ws.Range("K106").Value = 3994.5
$ws.Range("L106").Value = 1800
$ws.Range("M106").Value = -3363.5
$ws.Range("N106").Value = -3062

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2769.5
$ws.Range("I102").Value = 2670.8572
$ws.Range("K102").Value = 2670.8572
$ws.Range("M102").Value = -1048.8572

$ws.Range("H132").Value = 1380.9
$ws.Range("I132").Value = 1258.8572
$ws.Range("J132").Value = 1665.6666
$ws.Range("K132").Value = 3776.5716
$ws.Range("L132").Value = 4996.9998
$ws.Range("M132").Value = -1246.5716
$ws.Range("N132").Value = -10056.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 6667000
$ws.Range("I7").Value = 20000000
$ws.Range("K7").Value = 20000000
$ws.Range("M7").Value = -19999887

$ws.Range("H135").Value = 62998.332
$ws.Range("J135").Value = 62998.332
$ws.Range("L135").Value = 62998.332
$ws.Range("N135").Value = -73138.33199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2066
$ws.Range("I58").Value = 1500
$ws.Range("K58").Value = 1500
$ws.Range("M58").Value = -1297

$ws.Range("H134").Value = 5998.5
$ws.Range("I134").Value = 5998.5
$ws.Range("K134").Value = 17995.5
$ws.Range("M134").Value = -15460.5

$ws.Range("H136").Value = 2066
$ws.Range("I136").Value = 1500
$ws.Range("K136").Value = 4500
$ws.Range("M136").Value = -1950

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 66.09999999999999
$ws.Range("I2").Value = 26.23077
$ws.Range("K2").Value = 157.38462
$ws.Range("M2").Value = -44.38461999999998

$ws.Range("H5").Value = 1125
$ws.Range("I5").Value = 1166.6666
$ws.Range("K5").Value = 3499.9998
$ws.Range("M5").Value = -3387.9998

$ws.Range("H7").Value = 445.8889
$ws.Range("I7").Value = 97.14286
$ws.Range("K7").Value = 291.42858
$ws.Range("M7").Value = -179.42858

$ws.Range("H11").Value = 2777973.5
$ws.Range("I11").Value = 2941356
$ws.Range("K11").Value = 8824068
$ws.Range("M11").Value = -8823928

$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("N25").Value = 0
$ws.Range("L25").ClearContents()

$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("N30").Value = 0
$ws.Range("L30").ClearContents()

$ws.Range("H38").Value = 96.333336
$ws.Range("I38").Value = 90
$ws.Range("J38").Value = 99.5
$ws.Range("K38").Value = 270
$ws.Range("L38").Value = 298.5
$ws.Range("M38").Value = 77
$ws.Range("N38").Value = -992.5

$ws.Range("H80").Value = 5566.5557
$ws.Range("I80").Value = 4699.6665
$ws.Range("K80").Value = 14098.9995
$ws.Range("M80").Value = -13162.9995

$ws.Range("H83").Value = 5566.5557
$ws.Range("I83").Value = 4699.6665
$ws.Range("K83").Value = 42296.9985
$ws.Range("M83").Value = -37616.9985

$ws.Range("H92").Value = 2374.5
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 2374.5
$ws.Range("K92").Value = 0
$ws.Range("M92").Value = 7123.5
$ws.Range("N92").Value = -9619.5
$ws.Range("L92").ClearContents()

$ws.Range("H131").Value = 1774.75
$ws.Range("J131").Value = 2600
$ws.Range("L131").Value = 7800
$ws.Range("N131").Value = -17880

$ws.Range("H135").Value = 1125
$ws.Range("I135").Value = 1166.6666
$ws.Range("K135").Value = 10499.9994
$ws.Range("M135").Value = -7964.999400000001

$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("N138").Value = 0
$ws.Range("L138").ClearContents()
$ws.Range("M138").ClearContents()

$ws.Range("H139").Value = 2961.1667
$ws.Range("I139").Value = 1089
$ws.Range("J139").Value = 4833.3335
$ws.Range("K139").Value = 3267
$ws.Range("L139").Value = 14500.0005
$ws.Range("M139").Value = 1873
$ws.Range("N139").Value = -24780.0005

$ws.Range("H140").Value = 560.7778
$ws.Range("I140").Value = 560.7778
$ws.Range("K140").Value = 1682.3334
$ws.Range("M140").Value = 3497.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 2421017.5
$ws.Range("I3").Value = 2505221
$ws.Range("J3").Value = 2000000
$ws.Range("K3").Value = 2505221
$ws.Range("L3").Value = 2000000
$ws.Range("M3").Value = -2505105
$ws.Range("N3").Value = -2000232

$ws.Range("H21").Value = 2666666.8
$ws.Range("I21").Value = 2666666.8
$ws.Range("K21").Value = 2666666.8
$ws.Range("M21").Value = -2666493.8

$ws.Range("H30").Value = 2666666.8
$ws.Range("I30").Value = 2666666.8
$ws.Range("K30").Value = 2666666.8
$ws.Range("M30").Value = -2666561.8

$ws.Range("H122").Value = 800
$ws.Range("J122").Value = 800
$ws.Range("L122").Value = 2400
$ws.Range("N122").Value = -7300

$ws.Range("H126").Value = 5044.5557
$ws.Range("I126").Value = 5044.5557
$ws.Range("K126").Value = 15133.6671
$ws.Range("M126").Value = -12663.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3968.3333
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 53000
$ws.Range("J82").Value = 53000
$ws.Range("L82").Value = 53000
$ws.Range("N82").Value = -53766

$ws.Range("H85").Value = 53000
$ws.Range("J85").Value = 53000
$ws.Range("L85").Value = 53000
$ws.Range("N85").Value = -55652

$ws.Range("H132").Value = 7081.909
$ws.Range("I132").Value = 6650.143
$ws.Range("J132").Value = 7837.5
$ws.Range("K132").Value = 19950.429
$ws.Range("L132").Value = 23512.5
$ws.Range("M132").Value = -17420.429
$ws.Range("N132").Value = -28572.5
